# Adding Hapan Stats Part 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Append the new Hapan Consortium ("(HA)") ship rows (70-82).
#    Column-A values must be written in this exact order so the shared-string
#    table indices line up with the target workbook (new uniques are appended
#    in first-seen order).
# ---------------------------------------------------------------------------

# Row 70 - section header (name + type only, no stats)
$ws.Cells.Item(70,1).Value = "(HA) Ray"
$ws.Cells.Item(70,2).Value = "Frigate"
$ws.Cells.Item(70,1).Interior.Color = 65535

# Row 71 - (HA) Flare
$ws.Cells.Item(71,1).Value = "(HA) Flare"
$ws.Cells.Item(71,2).Value = "Frigate"
$ws.Cells.Item(71,3).Value = 200
$ws.Cells.Item(71,4).Value = 30
$ws.Cells.Item(71,5).Value = 0
$ws.Cells.Item(71,6).Value = 1
$ws.Cells.Item(71,7).Value = 25
$ws.Cells.Item(71,8).Value = 650
$ws.Cells.Item(71,9).Value = 650
$ws.Cells.Item(71,10).Value = 0
$ws.Cells.Item(71,11).Value = 3.5
$ws.Cells.Item(71,12).Value = 0.6
$ws.Cells.Item(71,13).Value = 1

# Row 72 - (HA) Beta
$ws.Cells.Item(72,1).Value = "(HA) Beta"
$ws.Cells.Item(72,2).Value = "Frigate"
$ws.Cells.Item(72,3).Value = 220
$ws.Cells.Item(72,4).Value = 40
$ws.Cells.Item(72,5).Value = 0
$ws.Cells.Item(72,6).Value = 2
$ws.Cells.Item(72,7).Value = 25
$ws.Cells.Item(72,8).Value = 1550
$ws.Cells.Item(72,9).Value = 1550
$ws.Cells.Item(72,10).Value = 0
$ws.Cells.Item(72,11).Value = 2
$ws.Cells.Item(72,12).Value = 0.6
$ws.Cells.Item(72,13).Value = 2

# Row 73 - (HA) Express
$ws.Cells.Item(73,1).Value = "(HA) Express"
$ws.Cells.Item(73,2).Value = "Frigate"
$ws.Cells.Item(73,3).Value = 200
$ws.Cells.Item(73,4).Value = 30
$ws.Cells.Item(73,5).Value = 0
$ws.Cells.Item(73,6).Value = 1
$ws.Cells.Item(73,7).Value = 25
$ws.Cells.Item(73,8).Value = 1550
$ws.Cells.Item(73,9).Value = 1550
$ws.Cells.Item(73,10).Value = 0
$ws.Cells.Item(73,11).Value = 2
$ws.Cells.Item(73,12).Value = 0.6
$ws.Cells.Item(73,13).Value = 1

# Row 74 - (HA) Stella
$ws.Cells.Item(74,1).Value = "(HA) Stella"
$ws.Cells.Item(74,2).Value = "Frigate"
$ws.Cells.Item(74,3).Value = 320
$ws.Cells.Item(74,4).Value = 60
$ws.Cells.Item(74,5).Value = 75
$ws.Cells.Item(74,6).Value = 4
$ws.Cells.Item(74,7).Value = 40
$ws.Cells.Item(74,8).Value = 1750
$ws.Cells.Item(74,9).Value = 1750
$ws.Cells.Item(74,10).Value = 133
$ws.Cells.Item(74,11).Value = 2
$ws.Cells.Item(74,12).Value = 0.6
$ws.Cells.Item(74,13).Value = 1

# Row 75 - (HA) Olanji
$ws.Cells.Item(75,1).Value = "(HA) Olanji"
$ws.Cells.Item(75,2).Value = "Frigate"
$ws.Cells.Item(75,3).Value = 650
$ws.Cells.Item(75,4).Value = 75
$ws.Cells.Item(75,5).Value = 50
$ws.Cells.Item(75,6).Value = 6
$ws.Cells.Item(75,7).Value = 60
$ws.Cells.Item(75,8).Value = 2800
$ws.Cells.Item(75,9).Value = 2800
$ws.Cells.Item(75,10).Value = 133
$ws.Cells.Item(75,11).Value = 2
$ws.Cells.Item(75,12).Value = 0.6
$ws.Cells.Item(75,13).Value = 3

# Row 76 - (HA) Charubah
$ws.Cells.Item(76,1).Value = "(HA) Charubah"
$ws.Cells.Item(76,2).Value = "Frigate"
$ws.Cells.Item(76,3).Value = 750
$ws.Cells.Item(76,4).Value = 85
$ws.Cells.Item(76,5).Value = 75
$ws.Cells.Item(76,6).Value = 4
$ws.Cells.Item(76,7).Value = 50
$ws.Cells.Item(76,8).Value = 1750
$ws.Cells.Item(76,9).Value = 1750
$ws.Cells.Item(76,10).Value = 24
$ws.Cells.Item(76,11).Value = 2
$ws.Cells.Item(76,12).Value = 0.6
$ws.Cells.Item(76,13).Value = 2

# Row 77 - (HA) Terephon
$ws.Cells.Item(77,1).Value = "(HA) Terephon"
$ws.Cells.Item(77,2).Value = "Frigate"
$ws.Cells.Item(77,3).Value = 850
$ws.Cells.Item(77,4).Value = 95
$ws.Cells.Item(77,5).Value = 85
$ws.Cells.Item(77,6).Value = 8
$ws.Cells.Item(77,7).Value = 60
$ws.Cells.Item(77,8).Value = 4550
$ws.Cells.Item(77,9).Value = 4550
$ws.Cells.Item(77,10).Value = 24
$ws.Cells.Item(77,11).Value = 2
$ws.Cells.Item(77,12).Value = 0.6
$ws.Cells.Item(77,13).Value = 3

# Row 79 - (HA) Corona (written before row 78 so the shared-string order
# for "(HA) Corona" etc. precedes "(HA) Magnetar", matching the target file)
$ws.Cells.Item(79,1).Value = "(HA) Corona"
$ws.Cells.Item(79,2).Value = "CapitalShip"
$ws.Cells.Item(79,3).Value = 1400
$ws.Cells.Item(79,4).Value = 110
$ws.Cells.Item(79,5).Value = 65
$ws.Cells.Item(79,6).Value = 1
$ws.Cells.Item(79,7).Value = 95
$ws.Cells.Item(79,8).Value = 4950
$ws.Cells.Item(79,9).Value = 4950
$ws.Cells.Item(79,10).Value = 56
$ws.Cells.Item(79,11).Value = 5.5
$ws.Cells.Item(79,12).Value = 0.65
$ws.Cells.Item(79,13).Value = 3

# Row 80 - (HA) Neutron
$ws.Cells.Item(80,1).Value = "(HA) Neutron"
$ws.Cells.Item(80,2).Value = "CapitalShip"
$ws.Cells.Item(80,3).Value = 900
$ws.Cells.Item(80,4).Value = 60
$ws.Cells.Item(80,5).Value = 110
$ws.Cells.Item(80,6).Value = 1
$ws.Cells.Item(80,7).Value = 85
$ws.Cells.Item(80,8).Value = 5450
$ws.Cells.Item(80,9).Value = 5450
$ws.Cells.Item(80,10).Value = 133
$ws.Cells.Item(80,11).Value = 5.5
$ws.Cells.Item(80,12).Value = 0.65
$ws.Cells.Item(80,13).Value = 3

# Row 81 - (HA) Pulsar
$ws.Cells.Item(81,1).Value = "(HA) Pulsar"
$ws.Cells.Item(81,2).Value = "CapitalShip"
$ws.Cells.Item(81,3).Value = 900
$ws.Cells.Item(81,4).Value = 60
$ws.Cells.Item(81,5).Value = 110
$ws.Cells.Item(81,6).Value = 1
$ws.Cells.Item(81,7).Value = 85
$ws.Cells.Item(81,8).Value = 5450
$ws.Cells.Item(81,9).Value = 5450
$ws.Cells.Item(81,10).Value = 133
$ws.Cells.Item(81,11).Value = 5.5
$ws.Cells.Item(81,12).Value = 0.65
$ws.Cells.Item(81,13).Value = 3

# Row 82 - (HA) Mist
$ws.Cells.Item(82,1).Value = "(HA) Mist"
$ws.Cells.Item(82,2).Value = "CapitalShip"
$ws.Cells.Item(82,3).Value = 1600
$ws.Cells.Item(82,4).Value = 180
$ws.Cells.Item(82,5).Value = 55
$ws.Cells.Item(82,6).Value = 1
$ws.Cells.Item(82,7).Value = 95
$ws.Cells.Item(82,8).Value = 6550
$ws.Cells.Item(82,9).Value = 6550
$ws.Cells.Item(82,10).Value = 56
$ws.Cells.Item(82,11).Value = 6
$ws.Cells.Item(82,12).Value = 0.65
$ws.Cells.Item(82,13).Value = 3

# Row 78 - (HA) Magnetar (written last among the new data rows so that its
# shared string lands after Corona/Neutron/Pulsar/Mist, as in the target)
$ws.Cells.Item(78,1).Value = "(HA) Magnetar"
$ws.Cells.Item(78,2).Value = "CapitalShip"
$ws.Cells.Item(78,3).Value = 1200
$ws.Cells.Item(78,4).Value = 100
$ws.Cells.Item(78,5).Value = 50
$ws.Cells.Item(78,6).Value = 1
$ws.Cells.Item(78,7).Value = 90
$ws.Cells.Item(78,8).Value = 4750
$ws.Cells.Item(78,9).Value = 4750
$ws.Cells.Item(78,10).Value = 56
$ws.Cells.Item(78,11).Value = 5.5
$ws.Cells.Item(78,12).Value = 0.65
$ws.Cells.Item(78,13).Value = 3

# ---------------------------------------------------------------------------
# 2) Rename the three old "(HC)" ships to "(HA)" - written last so their new
#    shared strings are appended at the very end of the table.
# ---------------------------------------------------------------------------
$ws.Cells.Item(51,1).Value = "(HA) BattleDragon"
$ws.Cells.Item(30,1).Value = "(HA) StarHome"
$ws.Cells.Item(43,1).Value = "(HA) NovaCruiser"

# ---------------------------------------------------------------------------
# 3) Resize columns B:M to fit their new contents.
# ---------------------------------------------------------------------------
$ws.Range("B1:M82").Columns.AutoFit()

# ---------------------------------------------------------------------------
# 4) Freeze the header rows (1-19) and scroll/select to match the saved view.
# ---------------------------------------------------------------------------
$ws.Range("A20").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 5) Apply the AutoFilter on column A, limited to the Hapan/Imperial Remnant
#    (and PA Secutor / PA Venator) ships - this hides all the other rows.
# ---------------------------------------------------------------------------
$filterVals = @(
  "(HA) BattleDragon",
  "(HA) Beta",
  "(HA) Charubah",
  "(HA) Corona",
  "(HA) Express",
  "(HA) Flare",
  "(HA) Magnetar",
  "(HA) Mist",
  "(HA) Neutron",
  "(HA) NovaCruiser",
  "(HA) Olanji",
  "(HA) Pulsar",
  "(HA) Ray",
  "(HA) StarHome",
  "(HA) Stella",
  "(HA) Terephon",
  "(IR) Allegiance",
  "(IR) Altor",
  "(IR) Carrack",
  "(IR) Dreadnaught",
  "(IR) EscortCarrier",
  "(IR) Immobilizer",
  "(IR) ISDI",
  "(IR) ISDII",
  "(IR) Katana",
  "(IR) Lancer",
  "(IR) MTC",
  "(IR) Pellaeon",
  "(IR) Strikecruiser",
  "(IR) Thrawn",
  "(IR) TorpedoSphere",
  "(IR) VSDI",
  "(IR) VSDII",
  "(IR) WorldDevastator",
  "(PA) Secutor",
  "(PA) Venator"
)
$ws.Range("A1:AO82").AutoFilter(1, $filterVals, 7)

# ---------------------------------------------------------------------------
# 6) Register the hidden _FilterDatabase defined name Excel creates when an
#    AutoFilter is applied.
# ---------------------------------------------------------------------------
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=information!`$A`$1:`$AO`$82")
$fdb.Visible = $false

# ---------------------------------------------------------------------------
# 7) Restore the final selection (bottom pane, cell K84) as in the target.
# ---------------------------------------------------------------------------
$ws.Range("K84").Select()
